$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin'
$ws.Range("D2").Value = '76.475.18'
$ws.Range("E2").Value = '  +0.62%  '

# Row 3: 'Ethereum'
$ws.Range("D3").Value = '2.939.56'
$ws.Range("E3").Value = '  +1.63%  '

# Row 4: 'TetherUSD'
$ws.Range("E4").Value = '  -0.07%  '

# Row 5: 'Solana'
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '198.76'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +1.04%  '

# Row 6: 'BNB'
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '595.26'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -0.98%  '

# Row 7: 'USDC'
$ws.Range("E7").Value = '  +0.04%  '

# Row 8: 'XRP'
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.551'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -0.78%  '

# Row 9: 'Dogecoin'
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.201'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +3.63%  '

# Row 10: 'LidoStakedEther'
$ws.Range("D10").Value = '2.938.81'
$ws.Range("E10").Value = '  +1.56%  '

# Row 11: 'Cardano'
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.441'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +10.17%  '

# Row 12: 'TRON'
$ws.Range("E12").Value = '  +0.43%  '

# Row 13: 'WrappedliquidstakedEther2.0' -> 'Toncoin'
$ws.Range("B13").Value = 'Toncoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '4.89'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  -0.96%  '

# Row 14: 'Toncoin' -> 'WrappedliquidstakedEther2.0'
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '3.479.32'
$ws.Range("E14").Value = '  +1.29%  '

# Row 15: 'Avalanche'
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '28.44'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +3.13%  '

# Row 16: 'WrappedBTC'
$ws.Range("D16").Value = '76.483.45'
$ws.Range("E16").Value = '  +0.77%  '

# Row 17: 'ShibaInu'
$ws.Range("E17").Value = '  -0.63%  '

# Row 18: 'WrappedEther'
$ws.Range("D18").Value = '2.927.10'
$ws.Range("E18").Value = '  +0.95%  '

# Row 19: 'Chainlink'
$ws.Range("E19").Value = '  +7.24%  '

# Row 20: 'Uniswap'
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '8.76'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -2.95%  '

# Row 21: 'BitcoinCash'
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '375.11'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -2.21%  '

# Row 22: 'Polkadot'
$ws.Range("E22").Value = '  +3.93%  '

# Row 23: 'SuiNetwork'
$ws.Range("E23").Value = '  -2.33%  '

# Row 24: 'Litecoin'
$ws.Range("E24").Value = '  -0.25%  '

# Row 25: 'Dai'
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +0.00%  '

# Row 26: 'WrappedeETH'
$ws.Range("D26").Value = '3.092.41'
$ws.Range("E26").Value = '  +1.59%  '

# Row 27: 'NEARProtocol'
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '4.27'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -0.11%  '

# Row 28: 'Aptos'
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '9.65'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  -2.23%  '

# Row 29: 'PEPE'
$ws.Range("E29").Value = '  -0.37%  '

# Row 30: 'Binance-PegBSC-USD'
$ws.Range("E30").Value = '  +0.23%  '

# Row 31: 'InternetComputer(DFINITY)'
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '8.32'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +5.99%  '

# Row 32: 'Fetch.AI'
$ws.Range("E32").Value = '  -3.00%  '

# Row 33: 'Bittensor'
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '497.08'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -3.09%  '

# Row 34: 'PancakeSwap'
$ws.Range("E34").Value = '  -0.51%  '

# Row 35: 'FirstDigitalUSD'
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  -0.08%  '

# Row 36: 'Monero'
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '165.39'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +0.23%  '

# Row 37: 'EthereumClassic'
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '20.17'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -0.35%  '

# Row 38: 'PolygonEcosystemToken'
$ws.Range("E38").Value = '  +12.98%  '

# Row 39: 'Cronos'
$ws.Range("E39").Value = '  +18.36%  '

# Row 40: 'WhiteBITCoin'
$ws.Range("E40").Value = '  +1.38%  '

# Row 41: 'Kaspa'
$ws.Range("E41").Value = '  -3.69%  '

# Row 43: 'Aave'
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '179.69'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -2.64%  '

# Row 44: 'RenderToken'
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '4.92'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  -3.06%  '

# Row 45: 'Stacks'
$ws.Range("E45").Value = '  -2.21%  '

# Row 46: 'OKB'
$ws.Range("E46").Value = '  -1.17%  '

# Row 47: 'ImmutableX' -> 'ARBITRUM'
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.595'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +1.27%  '

# Row 48: 'ARBITRUM' -> 'ImmutableX'
$ws.Range("B48").Value = 'ImmutableX'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '1.18'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -4.10%  '

# Row 49: 'dogwifhat'
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '2.34'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -2.21%  '

# Row 50: 'Filecoin'
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '3.86'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +2.05%  '

# Row 51: 'Mantle'
$ws.Range("E51").Value = '  -1.23%  '
Write-Output "Applied cryptos update"
